$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Caqui" log. In the
# source sheet this shows up as a fresh row inserted right after the
# header/first few rows (before the former row 20), which pushes every
# following record down by one and appends a former last row at the
# bottom (new dimension A1:T118).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly data point
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(20, 3).Value = 'Coquimbo'
$ws.Cells.Item(20, 4).Value = 45063
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 'Fruta'
$ws.Cells.Item(20, 7).Value = 100107
$ws.Cells.Item(20, 8).Value = 'Otros'
$ws.Cells.Item(20, 9).Value = 100107001
$ws.Cells.Item(20, 10).Value = 'Caqui'
$ws.Cells.Item(20, 11).Value = 'Mankaki'
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 14
$ws.Cells.Item(20, 14).Value = 330000
$ws.Cells.Item(20, 15).Value = 340000
$ws.Cells.Item(20, 16).Value = 335000
$ws.Cells.Item(20, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 744
$ws.Cells.Item(20, 20).Value = 450
